# Update the cryptocurrency price/volume table with freshly scraped values.
# Generated for commit: "Updated cryptos list on Mon Jun 24 04:45:40 UTC 2024 with GitHub Actions"

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of row -> (new Price text, new Volume(1h) text). A $null entry means
# that column is left unchanged for that row.
$updates = @(
    @{ Row = 2;  D = "62.891.27"; E = "  -2.30%  " },
    @{ Row = 3;  D = "3.403.83";  E = "  -3.20%  " },
    @{ Row = 4;  D = "1.00";      E = "  +0.00%  " },
    @{ Row = 5;  D = "574.38";    E = "  -3.03%  " },
    @{ Row = 6;  D = $null;       E = "  -6.06%  " },
    @{ Row = 7;  D = $null;       E = "  +0.02%  " },
    @{ Row = 8;  D = "3.404.29";  E = "  -3.16%  " },
    @{ Row = 9;  D = $null;       E = "  -2.78%  " },
    @{ Row = 10; D = "7.39";      E = "  -3.31%  " },
    @{ Row = 11; D = $null;       E = "  -3.14%  " },
    @{ Row = 12; D = "0.378";     E = "  -3.06%  " },
    @{ Row = 13; D = "3.978.80";  E = "  -3.30%  " },
    @{ Row = 14; D = $null;       E = "  -0.78%  " },
    @{ Row = 15; D = "3.400.19";  E = "  -3.25%  " },
    @{ Row = 16; D = $null;       E = "  -4.81%  " },
    @{ Row = 17; D = "62.885.92"; E = "  -2.30%  " },
    @{ Row = 18; D = "24.79";     E = "  -4.04%  " },
    @{ Row = 19; D = "9.52";      E = "  -4.93%  " },
    @{ Row = 20; D = $null;       E = "  -1.55%  " },
    @{ Row = 21; D = "13.16";     E = "  -3.03%  " },
    @{ Row = 22; D = "376.59";    E = "  -4.51%  " },
    @{ Row = 23; D = "0.558";     E = "  -3.22%  " },
    @{ Row = 24; D = "3.538.67";  E = "  -3.24%  " },
    @{ Row = 25; D = $null;       E = "  -0.17%  " },
    @{ Row = 26; D = "71.99";     E = $null },
    @{ Row = 27; D = $null;       E = "  -8.11%  " },
    @{ Row = 28; D = $null;       E = "  +0.04%  " },
    @{ Row = 29; D = $null;       E = "  -5.77%  " },
    @{ Row = 30; D = "2.15";      E = "  -4.81%  " },
    @{ Row = 31; D = "7.85";      E = "  -5.62%  " },
    @{ Row = 32; D = $null;       E = "  -4.90%  " },
    @{ Row = 33; D = "0.151";     E = "  -3.62%  " },
    @{ Row = 34; D = $null;       E = "  -0.02%  " },
    @{ Row = 35; D = "3.431.78";  E = "  -3.27%  " },
    @{ Row = 36; D = "22.73";     E = "  -2.99%  " },
    @{ Row = 37; D = "5.28";      E = "  -1.59%  " },
    @{ Row = 38; D = "164.98";    E = "  -1.18%  " },
    @{ Row = 39; D = "6.72";      E = "  -3.63%  " },
    @{ Row = 40; D = $null;       E = "  -4.46%  " },
    @{ Row = 41; D = "0.0756";    E = "  -4.40%  " },
    @{ Row = 42; D = "1.00";      E = "  -0.07%  " },
    @{ Row = 43; D = $null;       E = "  -4.66%  " },
    @{ Row = 44; D = "41.57";     E = "  -1.78%  " },
    @{ Row = 45; D = $null;       E = "  -4.05%  " },
    @{ Row = 46; D = "1.57";      E = "  -5.57%  " },
    @{ Row = 47; D = "22.97";     E = "  -10.76%  " },
    @{ Row = 48; D = $null;       E = "  -8.45%  " },
    @{ Row = 49; D = $null;       E = "  -2.22%  " },
    @{ Row = 50; D = "2.252.21";  E = "  -5.92%  " },
    @{ Row = 51; D = "0.853";     E = "  -5.19%  " }
)

foreach ($update in $updates) {
    $row = $update.Row

    if ($null -ne $update.D) {
        # Prefix with an apostrophe so Excel keeps the value as text instead
        # of silently re-parsing it into a number (e.g. "1.00" -> 1), then
        # restore the default "Normal" style so no stray formatting/xf index
        # is left behind on the cell.
        $dCell = $ws.Cells.Item($row, 4)
        $dCell.Value = "'" + $update.D
        $dCell.Style = "Normal"
    }

    if ($null -ne $update.E) {
        $eCell = $ws.Cells.Item($row, 5)
        $eCell.Value = $update.E
    }
}
